# Update the cryptos worksheet with the latest scraped price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.903.18"
$ws.Range("E2").Value = "  -1.32%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.400.54"
$ws.Range("E3").Value = "  -2.26%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.77"
$ws.Range("E5").Value = "  -0.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.89"
$ws.Range("E6").Value = "  -2.84%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.504"
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.158"
$ws.Range("E9").Value = "  +3.85%  "

# Row 10
$ws.Range("E10").Value = "  -1.30%  "

# Row 11
$ws.Range("E11").Value = "  -2.44%  "

# Row 12
$ws.Range("E12").Value = "  -2.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "67.827.97"
$ws.Range("E13").Value = "  -1.15%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000169"
$ws.Range("E14").Value = "  -0.33%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.76"
$ws.Range("E15").Value = "  -3.50%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.28"
$ws.Range("E16").Value = "  -4.47%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "330.02"
$ws.Range("E17").Value = "  -3.28%  "

# Row 18
$ws.Range("E18").Value = "  -3.82%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.75"
$ws.Range("E19").Value = "  -1.07%  "

# Row 20
$ws.Range("E20").Value = "  +0.02%  "

# Row 21
$ws.Range("E21").Value = "  -4.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.57"
$ws.Range("E22").Value = "  -1.96%  "

# Row 23
$ws.Range("E23").Value = "  -2.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.04"
$ws.Range("E24").Value = "  -1.86%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0₃0793"
$ws.Range("E25").Value = "  -2.98%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.04"
$ws.Range("E26").Value = "  -1.97%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "419.59"
$ws.Range("E28").Value = "  -3.66%  "

# Row 29
$ws.Range("E29").Value = "  -0.95%  "

# Row 30
$ws.Range("E30").Value = "  -1.76%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.83"
$ws.Range("E31").Value = "  +0.42%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.98"
$ws.Range("E32").Value = "  -0.30%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.62"
$ws.Range("E34").Value = "  -1.07%  "

# Row 35
$ws.Range("E35").Value = "  -4.00%  "

# Row 36
$ws.Range("E36").Value = "  -3.04%  "

# Row 37
$ws.Range("E37").Value = "  -5.18%  "

# Row 38
$ws.Range("E38").Value = "  -1.83%  "

# Row 39
$ws.Range("E39").Value = "  -4.83%  "

# Row 40
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "128.74"
$ws.Range("E40").Value = "  -3.18%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.28"
$ws.Range("E41").Value = "  -1.94%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.93"
$ws.Range("E42").Value = "  -7.40%  "

# Row 43
$ws.Range("E43").Value = "  -1.92%  "

# Row 45
$ws.Range("E45").Value = "  -1.11%  "

# Row 46
$ws.Range("E46").Value = "  +0.73%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.10"
$ws.Range("E47").Value = "  -0.86%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.35"
$ws.Range("E48").Value = "  -5.96%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.40"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0202"
$ws.Range("E50").Value = "  -6.60%  "

# Row 51
$ws.Range("E51").Value = "  -0.60%  "
